# Rename via pdf name instead of most recent
#
# Populates the Claims sheet with two claims whose patient names were
# resolved from the PDF (instead of defaulting to the "most recent" name),
# which in turn feeds the Summary sheet's Total Claims / Total Paid rollups.

$wb = $excel.ActiveWorkbook
$claims = $wb.Worksheets.Item("Claims")
$summary = $wb.Worksheets.Item("Summary")

$serviceDate = Get-Date -Year 2024 -Month 5 -Day 18 -Hour 0 -Minute 0 -Second 0

# Row 2: McGee, Test
$claims.Cells.Item(2, 1).Value = $serviceDate
$claims.Cells.Item(2, 2).Value = "McGee, Test"
$claims.Cells.Item(2, 3).Value = "5/1/24 - 5/18/24"
$claims.Cells.Item(2, 4).Value = 1300

# Row 3: Anna, Mary
$claims.Cells.Item(3, 1).Value = $serviceDate
$claims.Cells.Item(3, 2).Value = "Anna, Mary"
$claims.Cells.Item(3, 3).Value = "5/1/24 - 5/18/24"
$claims.Cells.Item(3, 4).Value = 400

# Leave the cursor where data entry on the Claims sheet last landed.
$claims.Range("D13").Select()

# Restore focus to the Summary tab (the sheet that was active before editing).
$summary.Activate()

$wb.Save()
